$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two placeholder neural-network rows are replaced: the first becomes
# "Deep Learning" (taking the second row's rating), the second becomes the
# freshly-run "kNN" model (taking the old kNN rating). The now-redundant
# third row is cleared out entirely.
$ws.Range("A4").Value = "SVM polynomial"
$ws.Range("A10").Value = "Deep Learning"
$ws.Range("B10").Value = 93.1
$ws.Range("A11").Value = "kNN"
$ws.Range("B11").Value = 94.2
$ws.Range("A12:B12").ClearContents()

# Insert new section title above the existing table
$ws.Range("A2").Value = "CARBON ONLY"

# Add the new Silicon isotope data section further down the sheet
$ws.Range("A18").Value = "CARBON and SILICON"
$ws.Range("A19").Value = "Model"
$ws.Range("B19").Value = "Rating (%)"
$ws.Range("A20").Value = "kNN"
$ws.Range("B20").Value = 94.7

# Match the view state captured after entering the new data
$ws.Range("C21").Select()
$ws.Application.ActiveWindow.ScrollRow = 5
